# Refresh the "Buying Opportunity" stock-screener sheet:
#  - update the ticker lists in columns B (long buildup... actually
#    "Buying Opportunity"), C (support Zone), D (long buildup),
#    E (Short buildup) and F (FII ENTERING) for existing rows 2-24
#  - append 13 new rows (25-37) of additional "support Zone" tickers,
#    extending the used range from A1:F24 to A1:F37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the sheet down to row 37, inheriting row 24's formatting
#     (keeps column A on the existing bold/centered/bordered style
#     instead of minting a brand-new style index) ---
$ws.Range("A24:F24").Copy()
$ws.Range("A25:F37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2 ---
$ws.Cells.Item(2,2).Value = "NSE:AGI"
$ws.Cells.Item(2,3).Value = "NSE:ACC"
$ws.Cells.Item(2,4).Value = "NSE:HDFCAMC"
$ws.Cells.Item(2,5).Value = "NSE:BAJAJFINSV"
$ws.Cells.Item(2,6).Value = "NSE:HDFCAMC"

# --- Row 3 ---
$ws.Cells.Item(3,2).Value = "NSE:AHL"
$ws.Cells.Item(3,3).Value = "NSE:ALKALI"
$ws.Cells.Item(3,5).Value = "NSE:BAJFINANCE"
$ws.Cells.Item(3,6).Value = "NSE:JKCEMENT"

# --- Row 4 ---
$ws.Cells.Item(4,2).Value = "NSE:HDFCAMC"
$ws.Cells.Item(4,3).Value = "NSE:ARCHIDPLY"
$ws.Cells.Item(4,5).Value = "NSE:DEEPAKNTR"
$ws.Cells.Item(4,6).ClearContents()

# --- Row 5 ---
$ws.Cells.Item(5,2).Value = "NSE:IFBIND"
$ws.Cells.Item(5,3).Value = "NSE:BAJAJHCARE"
$ws.Cells.Item(5,5).Value = "NSE:LALPATHLAB"
$ws.Cells.Item(5,6).ClearContents()

# --- Row 6 ---
$ws.Cells.Item(6,2).Value = "NSE:JAMNAAUTO"
$ws.Cells.Item(6,3).Value = "NSE:BANSWRAS"
$ws.Cells.Item(6,5).Value = "NSE:MARICO"
$ws.Cells.Item(6,6).ClearContents()

# --- Row 7 ---
$ws.Cells.Item(7,2).Value = "NSE:JKCEMENT"
$ws.Cells.Item(7,3).Value = "NSE:BIGBLOC"
$ws.Cells.Item(7,5).Value = "NSE:NTPC"
$ws.Cells.Item(7,6).ClearContents()

# --- Row 8 ---
$ws.Cells.Item(8,2).Value = "NSE:JKTYRE"
$ws.Cells.Item(8,3).Value = "NSE:BODALCHEM"
$ws.Cells.Item(8,5).Value = "NSE:PAGEIND"
$ws.Cells.Item(8,6).ClearContents()

# --- Row 9 ---
$ws.Cells.Item(9,2).Value = "NSE:JUBLINDS"
$ws.Cells.Item(9,3).Value = "NSE:CENTRUM"
$ws.Cells.Item(9,5).Value = "NSE:PFC"

# --- Row 10 ---
$ws.Cells.Item(10,2).Value = "NSE:KELLTONTEC"
$ws.Cells.Item(10,3).Value = "NSE:CHOICEIN"
$ws.Cells.Item(10,5).ClearContents()

# --- Row 11 ---
$ws.Cells.Item(11,2).Value = "NSE:KNRCON"
$ws.Cells.Item(11,3).Value = "NSE:DBL"
$ws.Cells.Item(11,5).ClearContents()

# --- Row 12 ---
$ws.Cells.Item(12,2).Value = "NSE:MAHLOG"
$ws.Cells.Item(12,3).Value = "NSE:DENORA"
$ws.Cells.Item(12,5).ClearContents()

# --- Row 13 ---
$ws.Cells.Item(13,2).Value = "NSE:NAM-INDIA"
$ws.Cells.Item(13,3).Value = "NSE:DISHTV"
$ws.Cells.Item(13,5).ClearContents()

# --- Row 14 ---
$ws.Cells.Item(14,2).Value = "NSE:PITTIENG"
$ws.Cells.Item(14,3).Value = "NSE:EXXARO"

# --- Row 15 ---
$ws.Cells.Item(15,2).Value = "NSE:POWERMECH"
$ws.Cells.Item(15,3).Value = "NSE:FACT"

# --- Row 16 ---
$ws.Cells.Item(16,2).ClearContents()
$ws.Cells.Item(16,3).Value = "NSE:HMVL"

# --- Row 17 ---
$ws.Cells.Item(17,2).ClearContents()
$ws.Cells.Item(17,3).Value = "NSE:IDBI"

# --- Row 18 ---
$ws.Cells.Item(18,2).ClearContents()
$ws.Cells.Item(18,3).Value = "NSE:INDIACEM"

# --- Row 19 ---
$ws.Cells.Item(19,2).ClearContents()
$ws.Cells.Item(19,3).Value = "NSE:INDIGOPNTS"

# --- Row 20 ---
$ws.Cells.Item(20,2).ClearContents()
$ws.Cells.Item(20,3).Value = "NSE:JPASSOCIAT"

# --- Row 21 ---
$ws.Cells.Item(21,2).ClearContents()
$ws.Cells.Item(21,3).Value = "NSE:KHAICHEM"

# --- Row 22 ---
$ws.Cells.Item(22,2).ClearContents()
$ws.Cells.Item(22,3).Value = "NSE:LYKALABS"

# --- Row 23 ---
$ws.Cells.Item(23,2).ClearContents()
$ws.Cells.Item(23,3).Value = "NSE:MAITHANALL"

# --- Row 24 ---
$ws.Cells.Item(24,2).ClearContents()
$ws.Cells.Item(24,3).Value = "NSE:MEDICAMEQ"

# --- New rows 25-37 (index in column A, ticker in column C) ---
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,3).Value = "NSE:MINDACORP"

$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,3).Value = "NSE:MOLDTECH"

$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,3).Value = "NSE:MOLDTKPAC"

$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,3).Value = "NSE:NBCC"

$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,3).Value = "NSE:NELCO"

$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,3).Value = "NSE:ORTINLAB"

$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,3).Value = "NSE:PALREDTEC"

$ws.Cells.Item(32,1).Value = 30
$ws.Cells.Item(32,3).Value = "NSE:PARACABLES"

$ws.Cells.Item(33,1).Value = 31
$ws.Cells.Item(33,3).Value = "NSE:PENINLAND"

$ws.Cells.Item(34,1).Value = 32
$ws.Cells.Item(34,3).Value = "NSE:PROZONER"

$ws.Cells.Item(35,1).Value = 33
$ws.Cells.Item(35,3).Value = "NSE:RAJSREESUG"

$ws.Cells.Item(36,1).Value = 34
$ws.Cells.Item(36,3).Value = "NSE:RAMASTEEL"

$ws.Cells.Item(37,1).Value = 35
$ws.Cells.Item(37,3).Value = "NSE:REPL"
